# Reorders the "Test Cases" table so the "description" column precedes the
# "uuid" column (mirroring an upstream Power Query field reorder), applies a
# number-format style to the boolean "passed" column, and restores the
# original column widths / selection for the swapped layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$lastRow = 10
$lastCol = 8

# --- Swap the data held in columns A (uuid) and B (description) ----------
for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aVal = $aCell.Value2
    $bVal = $bCell.Value2
    $aCell.Value2 = $bVal
    $bCell.Value2 = $aVal
}

# --- Swap the column widths to match the new column contents -------------
$widthA = $ws.Columns.Item(1).ColumnWidth
$widthB = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = $widthB
$ws.Columns.Item(2).ColumnWidth = $widthA

# --- Apply a number-format style to the "passed" boolean column ----------
$ws.Range("H2:H10").NumberFormat = "General"

# --- Update the active selection ------------------------------------------
$ws.Range("B4").Select()
